$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title timestamp
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 09:20"

# Row 8: Alemania -> Alemania
$ws.Range("B8").Value = 67051
$ws.Range("C8").Value = 166
$ws.Range("E8").Value = 52901
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 650

# Row 17: Austria -> Austria
$ws.Range("B17").Value = 9705
$ws.Range("C17").Value = 87
$ws.Range("E17").Value = 8961

# Row 22: Australia -> Australia
$ws.Range("D22").Value = 337
$ws.Range("E22").Value = 4201

# Row 25: Chequia -> Chequia
$ws.Range("B25").Value = 3002
$ws.Range("C25").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 24

# Row 37: Tailandia -> Tailandia
$ws.Range("D37").Value = 342
$ws.Range("E37").Value = 1299
$ws.Range("F37").Value = 23

# Row 66: Ucrania -> Ucrania
$ws.Range("B66").Value = 549
$ws.Range("C66").Value = 1
$ws.Range("E66").Value = 528

# Row 67: Lituania -> Lituania
$ws.Range("F67").Value = 27

# Row 69: Armenia -> Hungria
$ws.Range("A69").Value = "Hungria"
$ws.Range("B69").Value = 492
$ws.Range("C69").Value = 45
$ws.Range("D69").Value = 37
$ws.Range("E69").Value = 439
$ws.Range("F69").Value = 6
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 16

# Row 70: Hungria -> Armenia
$ws.Range("A70").Value = "Armenia"
$ws.Range("B70").Value = 482
$ws.Range("D70").Value = 30
$ws.Range("E70").Value = 449
$ws.Range("F70").Value = 15
$ws.Range("H70").Value = 3

# Row 72: Bulgaria -> Letonia
$ws.Range("A72").Value = "Letonia"
$ws.Range("B72").Value = 398
$ws.Range("C72").Value = 22
$ws.Range("D72").Value = 1
$ws.Range("E72").Value = 397
$ws.Range("F72").Value = 3
$ws.Range("H72").Value = 0

# Row 73: Bosnia y Herzegovina -> Bulgaria
$ws.Range("A73").Value = "Bulgaria"
$ws.Range("C73").Value = 20
$ws.Range("E73").Value = 354
$ws.Range("F73").Value = 13
$ws.Range("H73").Value = 8

# Row 74: Letonia -> Bosnia y Herzegovina
$ws.Range("A74").Value = "Bosnia y Herzegovina"
$ws.Range("B74").Value = 379
$ws.Range("C74").Value = 11
$ws.Range("D74").Value = 17
$ws.Range("E74").Value = 352
$ws.Range("F74").Value = 1
$ws.Range("H74").Value = 10

# Row 100: Malta -> Uzbekistan
$ws.Range("A100").Value = "Uzbekistan"
$ws.Range("B100").Value = 158
$ws.Range("C100").Value = 9
$ws.Range("D100").Value = 7
$ws.Range("E100").Value = 149
$ws.Range("F100").Value = 8
$ws.Range("H100").Value = 2

# Row 101: Ghana -> Malta
$ws.Range("A101").Value = "Malta"
$ws.Range("B101").Value = 156
$ws.Range("E101").Value = 154
$ws.Range("F101").Value = 4
$ws.Range("H101").Value = 0

# Row 102: Bielorrusia -> Ghana
$ws.Range("A102").Value = "Ghana"
$ws.Range("D102").Value = 2
$ws.Range("E102").Value = 145
$ws.Range("F102").Value = 1
$ws.Range("H102").Value = 5

# Row 103: Uzbekistan -> Bielorrusia
$ws.Range("A103").Value = "Bielorrusia"
$ws.Range("B103").Value = 152
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 47
$ws.Range("E103").Value = 105
$ws.Range("F103").Value = 2
$ws.Range("H103").Value = 0

# Row 109: Sri Lanka -> Sri Lanka
$ws.Range("D109").Value = 16
$ws.Range("E109").Value = 104

# Row 149: Tanzania -> Tanzania
$ws.Range("E149").Value = 17
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 1

# Row 160: Guinea Ecuatorial -> Dominica
$ws.Range("A160").Value = "Dominica"

# Row 161: Dominica -> Guinea Ecuatorial
$ws.Range("A161").Value = "Guinea Ecuatorial"

# Row 169: Suazilandia -> Granada
$ws.Range("A169").Value = "Granada"

# Row 170: Laos -> Suazilandia
$ws.Range("A170").Value = "Suazilandia"
$ws.Range("C170").Value = 0

# Row 171: Granada -> Laos
$ws.Range("A171").Value = "Laos"
$ws.Range("C171").Value = 1

# Row 173: Surinam -> Mozambique
$ws.Range("A173").Value = "Mozambique"

# Row 175: Guinea-Bisau -> Surinam
$ws.Range("A175").Value = "Surinam"

# Row 176: Mozambique -> Guinea-Bisau
$ws.Range("A176").Value = "Guinea-Bisau"

# Row 177: Guyana -> Zimbabue
$ws.Range("A177").Value = "Zimbabue"
$ws.Range("C177").Value = 1

# Row 178: Antigua y Barbuda -> Guyana
$ws.Range("A178").Value = "Guyana"
$ws.Range("B178").Value = 8
$ws.Range("H178").Value = 1

# Row 180: Zimbabue -> Antigua y Barbuda
$ws.Range("A180").Value = "Antigua y Barbuda"
$ws.Range("E180").Value = 7
$ws.Range("H180").Value = 0

# Row 191: Fiyi -> Montserrat
$ws.Range("A191").Value = "Montserrat"

# Row 192: Montserrat -> Fiyi
$ws.Range("A192").Value = "Fiyi"

# Row 198: Botsuana -> Liberia
$ws.Range("A198").Value = "Liberia"

# Row 199: Belice -> Islas Virgenes Britanicas
$ws.Range("A199").Value = "Islas Virgenes Britanicas"
$ws.Range("C199").Value = 1

# Row 200: Liberia -> Belice
$ws.Range("A200").Value = "Belice"

# Row 201: Islas Virgenes Britanicas -> Botsuana
$ws.Range("A201").Value = "Botsuana"
$ws.Range("C201").Value = 0

# Row 204: Papua Nueva Guinea -> Timor Oriental
$ws.Range("A204").Value = "Timor Oriental"

# Row 205: Timor Oriental -> Papua Nueva Guinea
$ws.Range("A205").Value = "Papua Nueva Guinea"
